# Updated cryptos list with GitHub Actions (price + 1h volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "61.718.81"
$ws.Range("E2").Value = "  -0.55%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "3.403.50"
$ws.Range("E3").Value = "  -0.34%  "
# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.24%  "
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.15"
$ws.Range("E5").Value = "  +0.80%  "
# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.23"
$ws.Range("E6").Value = "  +0.15%  "
# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -3.02%  "
# Row 8: USDC
$ws.Range("E8").Value = "  +0.01%  "
# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.722"
$ws.Range("E9").Value = "  -1.36%  "
# Row 10: Dogecoin
$ws.Range("E10").Value = "  -5.91%  "
# Row 11: Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.49"
$ws.Range("E11").Value = "  +0.02%  "
# Row 12: Polkadot
$ws.Range("E12").Value = "  +1.96%  "
# Row 13: ShibaInu
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -3.01%  "
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.938.25"
$ws.Range("E14").Value = "  -0.57%  "
# Row 15: TRON
$ws.Range("E15").Value = "  -0.05%  "
# Row 16: Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.32"
$ws.Range("E16").Value = "  -2.27%  "
# Row 17: WrappedEther
$ws.Range("D17").Value = "3.425.98"
$ws.Range("E17").Value = "  +0.60%  "
# Row 18: Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("E18").Value = "  +3.67%  "
# Row 19: Polygon
$ws.Range("E19").Value = "  +0.25%  "
# Row 20: WrappedBTC
$ws.Range("D20").Value = "61.743.42"
# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.03"
$ws.Range("E21").Value = "  +16.24%  "
# Row 22: Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.46"
$ws.Range("E22").Value = "  +1.52%  "
# Row 23: ImmutableX
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.26"
$ws.Range("E23").Value = "  +2.95%  "
# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.04"
$ws.Range("E24").Value = "  -0.05%  "
# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.29"
$ws.Range("E25").Value = "  +1.32%  "
# Row 26: Filecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.75"
$ws.Range("E26").Value = "  +10.60%  "
# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "32.96"
$ws.Range("E27").Value = "  -0.83%  "
# Row 28: LEO
$ws.Range("E28").Value = "  -0.55%  "
# Row 29: RenderToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.68"
$ws.Range("E29").Value = "  +1.25%  "
# Row 30: Toncoin
$ws.Range("E30").Value = "  -1.25%  "
# Row 31: Cosmos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.82"
$ws.Range("E31").Value = "  -0.88%  "
# Row 32: Kaspa
$ws.Range("E32").Value = "  -1.84%  "
# Row 34: InjectiveProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.80"
$ws.Range("E34").Value = "  -4.47%  "
# Row 35: Dai
$ws.Range("E35").Value = "  -0.69%  "
# Row 36: OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.78"
$ws.Range("E36").Value = "  +8.46%  "
# Row 37: VeChain
$ws.Range("E37").Value = "  -3.31%  "
# Row 38: FirstDigitalUSD
$ws.Range("E38").Value = "  +0.04%  "
# Row 39: Stacks
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.03"
$ws.Range("E39").Value = "  +4.37%  "
# Row 40: Monero
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.38"
$ws.Range("E40").Value = "  +5.81%  "
# Row 41: TheGraph
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.321"
$ws.Range("E41").Value = "  +2.90%  "
# Row 42: Stellar
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  -0.09%  "
# Row 43: LidoDAOToken
$ws.Range("E43").Value = "  -1.30%  "
# Row 44: ARBITRUM
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  +4.30%  "
# Row 45: WEMIXToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  +6.34%  "
# Row 46: NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.18"
$ws.Range("E46").Value = "  +1.85%  "
# Row 47: ThetaToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  +17.74%  "
# Row 48: Celestia
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.32"
$ws.Range("E48").Value = "  -1.69%  "
# Row 49: EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.14"
$ws.Range("E49").Value = "  +1.02%  "
# Row 50: PEPE
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0518"
$ws.Range("E50").Value = "  +16.09%  "
# Row 51: BitcoinSV
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.73"
$ws.Range("E51").Value = "  +14.22%  "
